$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 4 through 8 (Lasso, Support Vector Classifier, CART, Random Forest, XGBoost
# rows are removed entirely, leaving just Logistic Regression and LightGBM).
$ws.Range("A4:E8").EntireRow.Delete() | Out-Null

# Row 3 used to hold "Lasso" data; it now holds the "LightGBM" model (previously row 7).
$ws.Range("A3").Value = "LightGBM"

# Updated metric values (row 2 = Logistic Regression, row 3 = LightGBM).
$ws.Range("B2").Value = 0.8760330578512396
$ws.Range("C2").Value = 0.8766502597027046
$ws.Range("D2").Value = 0.8760330578512396
$ws.Range("E2").Value = 0.8747312434866046

$ws.Range("B3").Value = 0.8754820936639118
$ws.Range("C3").Value = 0.8763263685801709
$ws.Range("D3").Value = 0.8754820936639118
$ws.Range("E3").Value = 0.874127401684053
